$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new email address as a new row (A7), following the same pattern
# as the existing rows in the sheet.
$ws.Range("A7").Value = "regenerative.md@gmail.com"

# Give the new cell the same "Hyperlink" look the other email cells use.
$ws.Range("A7").Style = "Hyperlink"

# Wire up the mailto: hyperlink for the new cell (TextToDisplay matches the
# address, same as the source file's "display" attribute).
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:regenerative.md@gmail.com", "", "", "mailto:regenerative.md@gmail.com")

# Hyperlinks.Add() overwrites the cell text with the TextToDisplay value, and
# can touch the cell's formatting - restore both to match the source file.
$ws.Range("A7").Value = "regenerative.md@gmail.com"
$ws.Range("A7").Style = "Hyperlink"

# Match the recorded selection/cursor position after the edit.
$null = $ws.Range("A9").Select()

# Set the sheet to portrait orientation for printing.
$ws.PageSetup.Orientation = 1
